$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code sets")

# Insert a new row at 26 (duplicate of the Hypoxemia row) so the Hypoxemia /
# Hypoxia pair that used to live together in one PCORnet cell (row 25, F:
# "Hypoxemia`nHypoxia") gets split into two rows, each with a single term.
$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value2 = $ws.Cells.Item(25, 1).Value2   # A26 = Hypoxemia
$ws.Cells.Item(26, 2).Value2 = $ws.Cells.Item(25, 2).Value2   # B26 = Diagnosable condition
$ws.Cells.Item(26, 3).Value2 = $ws.Cells.Item(25, 3).Value2   # C26 = 1
$ws.Cells.Item(26, 4).Value2 = $ws.Cells.Item(25, 4).Value2   # D26 = Hypoxemia
$ws.Cells.Item(26, 6).Value  = "hypoxia"                       # F26 = hypoxia

# Clean up the PCORnet mapped-concepts column (F): drop the leading "*"
# markers and lower-case everything so it matches the PCORnet value-set
# naming convention; also split the merged multi-line values.
$ws.Range("F2").Value  = "difficulty concentrating/memory changes"
$ws.Range("F3").Value  = "difficulty concentrating/memory changes"
$ws.Range("F4").Value  = "headache"
$ws.Range("F5").Value  = "lightheadedness/fast heart rate"
$ws.Range("F6").Value  = "sleep disturbance"
$ws.Range("F7").Value  = "shortness of breath/cough"
$ws.Range("F8").Value  = "post-exertional malaise/persistent fatigue"
$ws.Range("F9").Value  = "post-exertional malaise/persistent fatigue"
$ws.Range("F10").Value = "smell and taste"
$ws.Range("F11").Value = "smell and taste"
$ws.Range("F12").Value = "bloating/constipation/diarrhea"
$ws.Range("F13").Value = "cardiovascular disease"
$ws.Range("F14").Value = "arrythmia"
$ws.Range("F15").Value = "blood clots"
$ws.Range("F16").Value = "chronic kidney disease"
$ws.Range("F19").Value = "diabetes"
$ws.Range("F21").Value = "migraine"
$ws.Range("F22").Value = "stroke"
$ws.Range("F23").Value = "mood disorders"
$ws.Range("F24").Value = "inserstitial lung disease"
$ws.Range("F25").Value = "hypoxemia"

$ws.Range("C1").Select()
